$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: in the source workbook the batsman name is "Pat Cummins" followed by a
# NON-BREAKING SPACE (U+00A0), not a regular space. Keep that exact character
# so the cell text matches the original data byte-for-byte.
$batsman = "Pat Cummins" + [char]0x00A0

# Full target table: header row + 11 data rows, columns A..K
$header = @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr")

$data = @(
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", $batsman, "17", "9", "1", "1", "188.88"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", $batsman, "12", "10", "1", "0", "120.00"),
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", $batsman, "53", "36", "5", "2", "147.22"),
    @(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kolkata Knight Riders", "Kings XI Punjab", $batsman, "5", "4", "0", "0", "125.00"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", $batsman, "15", "11", "0", "1", "136.36"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", $batsman, "5", "4", "1", "0", "125.00"),
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", $batsman, "1", "8", "0", "0", "12.50"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", $batsman, "4", "17", "0", "0", "23.52"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", $batsman, "0", "0", "0", "0", "-"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", $batsman, "33", "12", "1", "4", "275.00"),
    @(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", $batsman, "1", "3", "0", "0", "33.33")
)

# Pre-format the full used range as text so numeric-looking strings
# ("17", "120.00", "-", ...) stay stored as text, matching the source data.
$ws.Range("A1:K12").NumberFormat = "@"

for ($j = 0; $j -lt $header.Length; $j++) {
    $cell = $ws.Cells.Item(1, $j + 1)
    $cell.Value = $header[$j]
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $cell = $ws.Cells.Item($i + 2, $j + 1)
        $cell.Value = $row[$j]
    }
}
